$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 15874470
$ws.Range("I92").Value = 18519604
$ws.Range("J92").Value = 3666.6667
$ws.Range("K92").Value = 18519604
$ws.Range("L92").Value = 3666.6667
$ws.Range("M92").Value = -18518356
$ws.Range("N92").Value = -6162.6667

$ws.Range("H116").Value = 2999.75
$ws.Range("I116").Value = 2999.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2999.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 442.25
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 2054.709
$ws.Range("I132").Value = 1714.7142
$ws.Range("J132").Value = 4831.3335
$ws.Range("K132").Value = 5144.142599999999
$ws.Range("L132").Value = 14494.0005
$ws.Range("M132").Value = -2614.142599999999
$ws.Range("N132").Value = -19554.0005

$ws.Range("H141").Value = 4259
$ws.Range("I141").Value = 1772.125
$ws.Range("K141").Value = 5316.375
$ws.Range("M141").Value = -136.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 562312.1
$ws.Range("I32").Value = 629966.4
$ws.Range("J32").Value = 21078.5
$ws.Range("K32").Value = 629966.4
$ws.Range("L32").Value = 21078.5
$ws.Range("M32").Value = -629679.4
$ws.Range("N32").Value = -21652.5

$ws.Range("H61").Value = 2515.9138
$ws.Range("I61").Value = 2029.3889
$ws.Range("J61").Value = 3312.0454
$ws.Range("K61").Value = 2029.3889
$ws.Range("L61").Value = 3312.0454
$ws.Range("M61").Value = -1817.3889
$ws.Range("N61").Value = -3736.0454

$ws.Range("H74").Value = 2504.9546
$ws.Range("I74").Value = 2406.6316
$ws.Range("J74").Value = 3127.6667
$ws.Range("K74").Value = 2406.6316
$ws.Range("L74").Value = 3127.6667
$ws.Range("M74").Value = -1532.6316
$ws.Range("N74").Value = -4875.6667

$ws.Range("H77").Value = 2504.9546
$ws.Range("I77").Value = 2406.6316
$ws.Range("J77").Value = 3127.6667
$ws.Range("K77").Value = 12033.158
$ws.Range("L77").Value = 15638.3335
$ws.Range("M77").Value = -7665.158000000001
$ws.Range("N77").Value = -24374.3335

$ws.Range("H97").Value = 1172.9259
$ws.Range("I97").Value = 1067.65
$ws.Range("K97").Value = 1067.65
$ws.Range("M97").Value = -571.6500000000001

$ws.Range("H122").Value = 2786.8235
$ws.Range("I122").Value = 2282.7693
$ws.Range("J122").Value = 4425
$ws.Range("K122").Value = 6848.3079
$ws.Range("L122").Value = 13275
$ws.Range("M122").Value = -4398.3079
$ws.Range("N122").Value = -18175

$ws.Range("H136").Value = 2515.9138
$ws.Range("I136").Value = 2029.3889
$ws.Range("J136").Value = 3312.0454
$ws.Range("K136").Value = 6088.1667
$ws.Range("L136").Value = 9936.136200000001
$ws.Range("M136").Value = -3538.1667
$ws.Range("N136").Value = -15036.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1884.6364
$ws.Range("I86").Value = 1469.375
$ws.Range("J86").Value = 2992
$ws.Range("K86").Value = 1469.375
$ws.Range("L86").Value = 2992
$ws.Range("M86").Value = -346.375
$ws.Range("N86").Value = -5238

$ws.Range("H89").Value = 1884.6364
$ws.Range("I89").Value = 1469.375
$ws.Range("J89").Value = 2992
$ws.Range("K89").Value = 7346.875
$ws.Range("L89").Value = 14960
$ws.Range("M89").Value = -1730.875
$ws.Range("N89").Value = -26192

$ws.Range("H94").Value = 960
$ws.Range("I94").Value = 948.3889
$ws.Range("J94").Value = 974.9286
$ws.Range("K94").Value = 948.3889
$ws.Range("L94").Value = 974.9286
$ws.Range("M94").Value = -497.3889
$ws.Range("N94").Value = -1876.9286

$ws.Range("H107").Value = 1436.25
$ws.Range("I107").Value = 998
$ws.Range("K107").Value = 998
$ws.Range("M107").Value = 922

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8582.306
$ws.Range("I31").Value = 1787.6
$ws.Range("K31").Value = 1787.6
$ws.Range("M31").Value = -1492.6

$ws.Range("H34").Value = 8582.306
$ws.Range("I34").Value = 1787.6
$ws.Range("K34").Value = 1787.6
$ws.Range("M34").Value = -1585.6

$ws.Range("H134").Value = 2811.3125
$ws.Range("I134").Value = 2706.3147
$ws.Range("K134").Value = 8118.9441
$ws.Range("M134").Value = -5583.9441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 21819442
$ws.Range("I4").Value = 20000838
$ws.Range("J4").Value = 25002000
$ws.Range("K4").Value = 60002514
$ws.Range("L4").Value = 75006000
$ws.Range("M4").Value = -60002402
$ws.Range("N4").Value = -75006224

$ws.Range("H122").Value = 3065.9023
$ws.Range("I122").Value = 374.31818
$ws.Range("K122").Value = 3368.86362
$ws.Range("M122").Value = -918.8636200000001

$ws.Range("H126").Value = 1872.1111
$ws.Range("I126").Value = 1215
$ws.Range("J126").Value = 2059.8572
$ws.Range("K126").Value = 3645
$ws.Range("L126").Value = 6179.571599999999
$ws.Range("M126").Value = 1295
$ws.Range("N126").Value = -16059.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 1733.3334
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 5200.0002
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2750.0002
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 3343.318
$ws.Range("I132").Value = 3287.3794
$ws.Range("K132").Value = 9862.138199999999
$ws.Range("M132").Value = -7332.138199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3753.125
$ws.Range("I122").Value = 3275
$ws.Range("J122").Value = 3912.5
$ws.Range("K122").Value = 9825
$ws.Range("L122").Value = 11737.5
$ws.Range("M122").Value = -7375
$ws.Range("N122").Value = -16637.5

$ws.Range("H132").Value = 2736.2432
$ws.Range("I132").Value = 2700.9524
$ws.Range("J132").Value = 2782.5625
$ws.Range("K132").Value = 8102.8572
$ws.Range("L132").Value = 8347.6875
$ws.Range("M132").Value = -5572.8572
$ws.Range("N132").Value = -13407.6875

$ws.Range("H136").Value = 2646813.8
$ws.Range("I136").Value = 916.8049
$ws.Range("K136").Value = 2750.4147
$ws.Range("M136").Value = -200.4146999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 44464.5
$ws.Range("J46").Value = 44464.5
$ws.Range("L46").Value = 44464.5
$ws.Range("N46").Value = -44926.5

$ws.Range("H122").Value = 1873.4634
$ws.Range("I122").Value = 1890.1
$ws.Range("J122").Value = 1828.091
$ws.Range("K122").Value = 5670.299999999999
$ws.Range("L122").Value = 5484.272999999999
$ws.Range("M122").Value = -3220.299999999999
$ws.Range("N122").Value = -10384.273

$ws.Range("H132").Value = 2733983.8
$ws.Range("I132").Value = 1763.3422
$ws.Range("J132").Value = 7248087
$ws.Range("K132").Value = 5290.0266
$ws.Range("L132").Value = 21744261
$ws.Range("M132").Value = -2760.0266
$ws.Range("N132").Value = -21749321

$ws.Range("H134").Value = 44464.5
$ws.Range("J134").Value = 44464.5
$ws.Range("L134").Value = 133393.5
$ws.Range("N134").Value = -138463.5

$ws.Range("H136").Value = 2079.3965
$ws.Range("I136").Value = 1597.4082
$ws.Range("J136").Value = 4703.5557
$ws.Range("K136").Value = 4792.2246
$ws.Range("L136").Value = 14110.6671
$ws.Range("M136").Value = -2242.2246
$ws.Range("N136").Value = -19210.6671
